$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Valor Mora total
$ws.Range("E11").Value = 228800

# Update Periodo Mora values (shift the period list forward by one: drop 2504, add 2508)
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2507"
$ws.Range("E19").Value = "2508"

# Update Valor Mora for the last worker row to match the others
$ws.Range("F19").Value = 57200
